$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the "last updated" timestamp string (row 1 / cell A1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Julio de 2020 a las 16:28"

# ---------------------------------------------------------------------------
# 2) Countries whose case counters were refreshed in place (rank unchanged).
# ---------------------------------------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 3294539
$ws.Range("C4").Value = 2753
$ws.Range("D4").Value = 1460771
$ws.Range("E4").Value = 1697033
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 64
$ws.Range("H4").Value = 136735

# Row 6: India
$ws.Range("B6").Value = 830763
$ws.Range("C6").Value = 8160
$ws.Range("D6").Value = 522631
$ws.Range("E6").Value = 285877
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 111
$ws.Range("H6").Value = 22255

# Row 19: Alemania
$ws.Range("B19").Value = 199652
$ws.Range("C19").Value = 64
$ws.Range("D19").Value = 184500
$ws.Range("E19").Value = 6020
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 2
$ws.Range("H19").Value = 9132

# Row 75: Australia
$ws.Range("B75").Value = 9553
$ws.Range("C75").Value = 194
$ws.Range("D75").Value = 7730
$ws.Range("E75").Value = 1716
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 107

# Row 77: Noruega
$ws.Range("B77").Value = 8976
$ws.Range("C77").Value = 2
$ws.Range("D77").Value = 8138
$ws.Range("E77").Value = 586
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 252

# Row 89: Tayikistan
$ws.Range("B89").Value = 6506
$ws.Range("C89").Value = 49
$ws.Range("D89").Value = 5176
$ws.Range("E89").Value = 1275
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 55

# Row 92: Estado de Palestina
$ws.Range("B92").Value = 5931
$ws.Range("C92").Value = 380
$ws.Range("D92").Value = 536
$ws.Range("E92").Value = 5364
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 4
$ws.Range("H92").Value = 31

# Row 110: Sri Lanka
$ws.Range("B110").Value = 2459
$ws.Range("C110").Value = 5
$ws.Range("D110").Value = 1980
$ws.Range("E110").Value = 468
$ws.Range("F110").Value = 0
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 11

# Row 111: Cuba
$ws.Range("B111").Value = 2420
$ws.Range("C111").Value = 7
$ws.Range("D111").Value = 2254
$ws.Range("E111").Value = 79
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 87

# Row 164: Birmania
$ws.Range("B164").Value = 330
$ws.Range("C164").Value = 4
$ws.Range("D164").Value = 260
$ws.Range("E164").Value = 64
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 6

# ---------------------------------------------------------------------------
# 3) Countries that swapped ranking with a neighbour: the newly higher-ranked
#    country takes brand new totals, the displaced one slides down a row with
#    the figures the displacing country used to show.
# ---------------------------------------------------------------------------
# Suecia(28)/Indonesia(29)/Irak(30)  ->  Irak(28)/Suecia(29)/Indonesia(30)
$ws.Range("A30").Value = "Indonesia"
$ws.Range("B30").Value = 74018
$ws.Range("C30").Value = 1671
$ws.Range("D30").Value = 34719
$ws.Range("E30").Value = 35764
$ws.Range("F30").Value = 0
$ws.Range("G30").Value = 66
$ws.Range("H30").Value = 3535

$ws.Range("A29").Value = "Suecia"
$ws.Range("B29").Value = 74898
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 5526

$ws.Range("A28").Value = "Irak"
$ws.Range("B28").Value = 75194
$ws.Range("C28").Value = 2734
$ws.Range("D28").Value = 43079
$ws.Range("E28").Value = 29060
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 95
$ws.Range("H28").Value = 3055

# Singapur(41)/Portugal(42)  ->  Portugal(41)/Singapur(42)
$ws.Range("A42").Value = "Singapur"
$ws.Range("B42").Value = 45783
$ws.Range("C42").Value = 170
$ws.Range("D42").Value = 41780
$ws.Range("E42").Value = 3977
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 26

$ws.Range("A41").Value = "Portugal"
$ws.Range("B41").Value = 46221
$ws.Range("C41").Value = 542
$ws.Range("D41").Value = 30655
$ws.Range("E41").Value = 13912
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 1654

# Yemen(128)/Libia(129)  ->  Libia(128)/Yemen(129)
$ws.Range("A129").Value = "Yemen"
$ws.Range("B129").Value = 1380
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 630
$ws.Range("E129").Value = 386
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 364

$ws.Range("A128").Value = "Libia"
$ws.Range("B128").Value = 1389
$ws.Range("C128").Value = 47
$ws.Range("D128").Value = 340
$ws.Range("E128").Value = 1011
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 38

